$wb = $excel.ActiveWorkbook

# --- Rename header in "Weekly Quantity" sheet ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- Rename header in "Monthly Trend" sheet ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet at the end ---
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "PO Forecast"

# Copy formatting (bold/bordered/centered header style, date-number-format column)
# from the "Weekly Quantity" sheet so the new sheet matches the workbook's look & feel.
$wsWeekly.Range("A1:B1").Copy($newSheet.Range("A1:B1"))
$wsWeekly.Range("A1:B1").Copy($newSheet.Range("C1:D1"))
$wsWeekly.Range("A2:A18").Copy($newSheet.Range("A2:A18"))
$wsWeekly.Range("A2:A9").Copy($newSheet.Range("A19:A26"))

# Headers (overwrite copied text with the new sheet's own header labels)
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Data rows
$newSheet.Range("A2").Value = 45354.99999999999
$newSheet.Range("B2").Value = 29
$newSheet.Range("C2").Value = -0.583532548516661
$newSheet.Range("D2").Value = 60.39808002107727
$newSheet.Range("A3").Value = 45368.99999999999
$newSheet.Range("B3").Value = 29
$newSheet.Range("C3").Value = -4.280042988910019
$newSheet.Range("D3").Value = 59.4158074461247
$newSheet.Range("A4").Value = 45382.99999999999
$newSheet.Range("B4").Value = 28
$newSheet.Range("C4").Value = -3.601169056456649
$newSheet.Range("D4").Value = 58.16125737934308
$newSheet.Range("A5").Value = 45445.99999999999
$newSheet.Range("B5").Value = 25
$newSheet.Range("C5").Value = -6.622774087277858
$newSheet.Range("D5").Value = 55.31281594840316
$newSheet.Range("A6").Value = 45452.99999999999
$newSheet.Range("B6").Value = 24
$newSheet.Range("C6").Value = -7.579724689094524
$newSheet.Range("D6").Value = 56.13669337763232
$newSheet.Range("A7").Value = 45459.99999999999
$newSheet.Range("B7").Value = 24
$newSheet.Range("C7").Value = -8.076317163087984
$newSheet.Range("D7").Value = 56.51039312247143
$newSheet.Range("A8").Value = 45466.99999999999
$newSheet.Range("B8").Value = 24
$newSheet.Range("C8").Value = -7.559135052188438
$newSheet.Range("D8").Value = 58.83357817822088
$newSheet.Range("A9").Value = 45480.99999999999
$newSheet.Range("B9").Value = 23
$newSheet.Range("C9").Value = -8.711095591067284
$newSheet.Range("D9").Value = 55.18124805540552
$newSheet.Range("A10").Value = 45487.99999999999
$newSheet.Range("B10").Value = 23
$newSheet.Range("C10").Value = -10.41193809560292
$newSheet.Range("D10").Value = 55.01795177315501
$newSheet.Range("A11").Value = 45515.99999999999
$newSheet.Range("B11").Value = 21
$newSheet.Range("C11").Value = -11.35184674222124
$newSheet.Range("D11").Value = 54.06840225717007
$newSheet.Range("A12").Value = 45529.99999999999
$newSheet.Range("B12").Value = 20
$newSheet.Range("C12").Value = -9.557582177449991
$newSheet.Range("D12").Value = 50.49022557405488
$newSheet.Range("A13").Value = 45536.99999999999
$newSheet.Range("B13").Value = 20
$newSheet.Range("C13").Value = -10.46731855333987
$newSheet.Range("D13").Value = 51.81293016480551
$newSheet.Range("A14").Value = 45543.99999999999
$newSheet.Range("B14").Value = 20
$newSheet.Range("C14").Value = -12.06384073319787
$newSheet.Range("D14").Value = 50.74169927075113
$newSheet.Range("A15").Value = 45578.99999999999
$newSheet.Range("B15").Value = 18
$newSheet.Range("C15").Value = -12.98498983903117
$newSheet.Range("D15").Value = 50.31509994917202
$newSheet.Range("A16").Value = 45585.99999999999
$newSheet.Range("B16").Value = 18
$newSheet.Range("C16").Value = -13.75760578937073
$newSheet.Range("D16").Value = 49.06197453175329
$newSheet.Range("A17").Value = 45592.99999999999
$newSheet.Range("B17").Value = 17
$newSheet.Range("C17").Value = -15.95306522777737
$newSheet.Range("D17").Value = 48.73979446631426
$newSheet.Range("A18").Value = 45599.99999999999
$newSheet.Range("B18").Value = 17
$newSheet.Range("C18").Value = -13.21259815164336
$newSheet.Range("D18").Value = 48.01466004912808
$newSheet.Range("A19").Value = 45606.99999999999
$newSheet.Range("B19").Value = 16
$newSheet.Range("C19").Value = -17.42965551239562
$newSheet.Range("D19").Value = 48.62625686708801
$newSheet.Range("A20").Value = 45613.99999999999
$newSheet.Range("B20").Value = 16
$newSheet.Range("C20").Value = -16.46360787334077
$newSheet.Range("D20").Value = 47.04657264631606
$newSheet.Range("A21").Value = 45620.99999999999
$newSheet.Range("B21").Value = 16
$newSheet.Range("C21").Value = -14.09261616713319
$newSheet.Range("D21").Value = 48.80298101110177
$newSheet.Range("A22").Value = 45627.99999999999
$newSheet.Range("B22").Value = 15
$newSheet.Range("C22").Value = -13.60274953708976
$newSheet.Range("D22").Value = 47.22421606859573
$newSheet.Range("A23").Value = 45634.99999999999
$newSheet.Range("B23").Value = 15
$newSheet.Range("C23").Value = -14.15762907978234
$newSheet.Range("D23").Value = 46.65915365039277
$newSheet.Range("A24").Value = 45641.99999999999
$newSheet.Range("B24").Value = 15
$newSheet.Range("C24").Value = -15.68139181361605
$newSheet.Range("D24").Value = 47.17599439565127
$newSheet.Range("A25").Value = 45648.99999999999
$newSheet.Range("B25").Value = 14
$newSheet.Range("C25").Value = -18.83821581087807
$newSheet.Range("D25").Value = 47.44984371584771
$newSheet.Range("A26").Value = 45655.99999999999
$newSheet.Range("B26").Value = 14
$newSheet.Range("C26").Value = -17.81030704619955
$newSheet.Range("D26").Value = 43.31242468309458
Write-Output "Edit applied successfully"
